$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F rows explaining the numeric encodings used for the categorical
# variables (sex, race, c_charge_degree, r/vr_charge_degree).
$ws.Range("F7").Value = "índices das variáveis categóricas"
$ws.Range("F8").Value = "sex = {'Female':1, 'Male':2}"
$ws.Range("F9").Value = "race = {'African-American':1, 'Asian':2, 'Caucasian':3, 'Hispanic':4, 'Native American':5, 'Other':6}"
$ws.Range("F10").Value = "c_charge_degree = {'F':1, 'M':2}"
$ws.Range("F11").Value = "r_vr_charge_degree = {np.nan:0, '(CO3)':1, '(F1)':2, '(F2)':3, '(F3)':4, '(F5)':5, '(F6)':6, '(F7)':7, '(M1)':8, '(M2)':9, '(MO3)':10}"

# Widen column F so the long descriptions fit (mirrors the bestFit resize
# Excel performs automatically when such text is entered).
$ws.Columns.Item(6).ColumnWidth = 108.33
